$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.424.35"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.945.65"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.54"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.88"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("D7").Value = "3.943.62"
$ws.Range("E7").Value = "  +2.28%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.46"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("E13").Value = "  +4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.84"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").Value = "4.601.57"
$ws.Range("E15").Value = "  +2.59%  "
$ws.Range("D16").Value = "3.969.00"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("D17").Value = "69.287.14"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.51"
$ws.Range("E18").Value = "  -0.89%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.41"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.02"
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "499.08"
$ws.Range("E22").Value = "  +1.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.735"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000167"
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.31"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.21"
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.25"
$ws.Range("E28").Value = "  +1.63%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "4.097.69"
$ws.Range("E31").Value = "  +2.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.41"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.81"
$ws.Range("E33").Value = "  -2.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.19"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("D35").Value = "3.917.39"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.04"
$ws.Range("E37").Value = "  +0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.05"
$ws.Range("E38").Value = "  +1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.30"
$ws.Range("E40").Value = "  +8.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.995"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.323"
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.03"
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "433.89"
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.18"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.61"
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000283"
$ws.Range("E48").Value = "  +23.79%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0364"
$ws.Range("E49").Value = "  +2.19%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.48"
$ws.Range("E50").Value = "  -0.54%  "
$ws.Range("D51").Value = "2.808.55"
$ws.Range("E51").Value = "  -1.38%  "
